$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Row 27: fill previously-blank cells with the literal string "nan"
$ws.Range("B27").Value = "nan"
$ws.Range("C27").Value = "nan"
$ws.Range("D27").Value = "nan"
$ws.Range("E27").Value = "nan"
$ws.Range("F27").Value = "nan"
$ws.Range("G27").Value = "nan"
$ws.Range("H27").Value = "nan"
$ws.Range("I27").Value = "nan"
$ws.Range("J27").Value = "nan"
$ws.Range("K27").Value = "nan"
$ws.Range("N27").Value = "nan"
$ws.Range("Q27").Value = "nan"

# Row 28: new event row duplicated from row 27 (card id, date, event, correction, serviced by)
$ws.Range("A28").Value = "'18"
$ws.Range("A28").Style = "Normal"

$ws.Range("L28").Value = "'12/1/2026"
$ws.Range("L28").Style = "Normal"

$ws.Range("M28").Value = "صوت وفايبريشن عالي  ف مجموعه دليفري وبعد  معاينه وجد تاكل ف طاره عصاره خروج شريط"
$ws.Range("O28").Value = "تم تركيب الطاره عكس اتجاها لحين تصنيع طاره أخري وتم تنظيف رولين سير700"
$ws.Range("P28").Value = "م.محمد عبدالله ،محمود ايهاب،،سلامه،ابراهيم،مصطفي"
